$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.960.22"
$ws.Range("E2").Value = "  -2.76%  "
$ws.Range("D3").Value = "3.197.27"
$ws.Range("E3").Value = "  -1.53%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'600.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").Value = "'153.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.39%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").Value = "3.196.12"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("E9").Value = "  -3.47%  "
$ws.Range("E10").Value = "  -3.41%  "
$ws.Range("E11").Value = "  -3.03%  "
$ws.Range("E12").Value = "  -6.35%  "
$ws.Range("D13").Value = "'0.0000258"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.24%  "
$ws.Range("D14").Value = "'37.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.54%  "
$ws.Range("D15").Value = "3.732.23"
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").Value = "65.158.63"
$ws.Range("E16").Value = "  -2.60%  "
$ws.Range("D17").Value = "3.206.87"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "'7.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.11%  "
$ws.Range("D20").Value = "'485.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.52%  "
$ws.Range("D21").Value = "'14.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.21%  "
$ws.Range("D22").Value = "'0.722"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.24%  "
$ws.Range("D23").Value = "'7.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.25%  "
$ws.Range("D24").Value = "'14.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.74%  "
$ws.Range("D25").Value = "'84.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").Value = "'8.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.98%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "'0.129"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +36.53%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'2.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.90%  "
$ws.Range("E31").Value = "  -4.80%  "
$ws.Range("D32").Value = "'2.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.32%  "
$ws.Range("D33").Value = "'27.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.45%  "
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("E35").Value = "  -6.35%  "
$ws.Range("D36").Value = "'6.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.89%  "
$ws.Range("D37").Value = "'54.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("E38").Value = "  +9.21%  "
$ws.Range("D39").Value = "'476.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.71%  "
$ws.Range("D40").Value = "0.0₃0742"
$ws.Range("E40").Value = "  -2.94%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0407"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.39%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.126"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.00%  "
$ws.Range("D43").Value = "'8.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.63%  "
$ws.Range("D44").Value = "2.935.71"
$ws.Range("E44").Value = "  +2.24%  "
$ws.Range("D45").Value = "'2.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("D46").Value = "'0.282"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.08%  "
$ws.Range("D47").Value = "'27.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.76%  "
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D51").Value = "'120.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.16%  "
